$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update total "Valor Mora" (E11)
$ws.Cells.Item(11, 5).Value = 1383198

# 2) Insert 5 additional rows before the old closing (bordered) row 37, pushing
#    the old row 37 (and the blank rows + summary block after it) down to
#    make room for rows 37-41, while row 42 becomes the new closing row.
$ws.Range("B37:B41").EntireRow.Insert(-4121)

# Copy the formatting (styles) of row 36 (a normal interior row of the table)
# across the 5 freshly inserted rows 37-41 so they look like the rest of the table.
$ws.Range("B36:J36").Copy($ws.Range("B37:J41"))

# Clear the old text content of the worker/period columns so that the shared
# string table drops the now-unused strings and reassigns ids in the order
# the new values are written below (matches how the source workbook was
# rebuilt).
$ws.Range("C16:E42").ClearContents()

# 3) Re-populate the whole worker/period table (rows 16-42) with the new data.
#    Columns: B=Tipo Doc (CC), C=N Doc Trabajador, D=Nombre Trabajador,
#             E=Periodo Mora, F=Valor Mora, G=Salario Basico (unchanged 1300000)

$rows = @(
  @{R=16; C="33222782";    D="VICKY URANGO BELEÑO";      E="2505"; F=45066},
  @{R=17; C="33222782";    D="VICKY URANGO BELEÑO";      E="2504"; F=52000},
  @{R=18; C="33222782";    D="VICKY URANGO BELEÑO";      E="2503"; F=52000},
  @{R=19; C="33222782";    D="VICKY URANGO BELEÑO";      E="2502"; F=52000},
  @{R=20; C="33222782";    D="VICKY URANGO BELEÑO";      E="2501"; F=52000},
  @{R=21; C="33222782";    D="VICKY URANGO BELEÑO";      E="2412"; F=52000},
  @{R=22; C="33222782";    D="VICKY URANGO BELEÑO";      E="2411"; F=52000},
  @{R=23; C="1002241292";  D="OMAR TORRES SALLAS";       E="2504"; F=52000},
  @{R=24; C="1002241292";  D="OMAR TORRES SALLAS";       E="2503"; F=52000},
  @{R=25; C="1002241292";  D="OMAR TORRES SALLAS";       E="2502"; F=52000},
  @{R=26; C="1002241292";  D="OMAR TORRES SALLAS";       E="2501"; F=52000},
  @{R=27; C="1002241292";  D="OMAR TORRES SALLAS";       E="2412"; F=52000},
  @{R=28; C="1002241292";  D="OMAR TORRES SALLAS";       E="2411"; F=52000},
  @{R=29; C="1047471137";  D="CENIA LUZ HERRERA GODOY";  E="2505"; F=45066},
  @{R=30; C="1047471137";  D="CENIA LUZ HERRERA GODOY";  E="2504"; F=52000},
  @{R=31; C="1047471137";  D="CENIA LUZ HERRERA GODOY";  E="2503"; F=52000},
  @{R=32; C="1047471137";  D="CENIA LUZ HERRERA GODOY";  E="2502"; F=52000},
  @{R=33; C="1047471137";  D="CENIA LUZ HERRERA GODOY";  E="2501"; F=52000},
  @{R=34; C="1047471137";  D="CENIA LUZ HERRERA GODOY";  E="2412"; F=52000},
  @{R=35; C="1047471137";  D="CENIA LUZ HERRERA GODOY";  E="2411"; F=52000},
  @{R=36; C="1047470457";  D="INDIRA HERRERA GODOY";     E="2505"; F=45066},
  @{R=37; C="1047470457";  D="INDIRA HERRERA GODOY";     E="2504"; F=52000},
  @{R=38; C="1047470457";  D="INDIRA HERRERA GODOY";     E="2503"; F=52000},
  @{R=39; C="1047470457";  D="INDIRA HERRERA GODOY";     E="2502"; F=52000},
  @{R=40; C="1047470457";  D="INDIRA HERRERA GODOY";     E="2501"; F=52000},
  @{R=41; C="1047470457";  D="INDIRA HERRERA GODOY";     E="2412"; F=52000},
  @{R=42; C="1047470457";  D="INDIRA HERRERA GODOY";     E="2411"; F=52000}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = 1300000
}

Write-Host "Edit complete"
